$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set / correct the "should_be_included" (column I) AI relevance-check values
$ws.Range("I8").Value = "Yes"
$ws.Range("I11").Value = "Unsure"
$ws.Range("I12").Value = "No"
$ws.Range("I13").Value = "Unsure"
$ws.Range("I14").Value = "Yes"
$ws.Range("I15").Value = "No"
$ws.Range("I16").Value = "No"
$ws.Range("I17").Value = "Unsure"
$ws.Range("I18").Value = "No"
$ws.Range("I19").Value = "Unsure"
$ws.Range("I20").Value = "No"
$ws.Range("I21").Value = "Unsure"
$ws.Range("I22").Value = "No"

# Update the view state (scroll position, zoom, selection) to match author's last save
$ws.Application.ActiveWindow.Zoom = 70
$ws.Range("T39").Select()
$ws.Application.ActiveWindow.ScrollRow = 39
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Application.ActiveWindow.TopLeftCell = $ws.Range("A39")
